# Generate Report for Handback
# Rename the two tracked files and refresh their handoff/handback timestamps
# across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldFile1 = "67c54cca-20e7-401f-a0e0-c8583438eec2.md"
$newFile1 = "dc50655f-ca93-4596-8aeb-64a92162e418.md"
$oldFile2 = "a97d0c6a-328d-4e75-a747-e0be0fa92585.md"
$newFile2 = "ffffbbf77c52-2624-4515-b61d-df71def2d5ad.md"

$newDate1 = "2016-09-07 07:28:18"

$newHandoffZhCn = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.zh-cn.xlf"
$newHandoffDateZhCn = "2016-09-07 07:28:11"
$newHandbackDateZhCn = "2016-09-07 07:28:37"

$newHandoffDeDe = "dc50655f-ca93-4596-8aeb-64a92162e418.cdb0fd07554f79e3abeb17b25c507fe8dc4f629e.de-de.xlf"
$newHandbackDateDeDe = "2016-09-07 07:28:45"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("G2").Value = $newDate1

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("G3").Value = $newDate1

# Hyperlinks: keep the same target addresses / relationship ids (the
# generator does not touch the handoff-commit URL on rename), only the
# displayed text needs to change. Updating a Hyperlink object's properties
# in place appends a duplicate entry, so rebuild the collection instead.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/$oldFile1", [System.Type]::Missing, [System.Type]::Missing, "e2e\$newFile1") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/$oldFile2", [System.Type]::Missing, [System.Type]::Missing, "e2e\$newFile2") | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = $newHandoffZhCn
$wsZhCn.Range("H2").Value = $newHandoffDateZhCn
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = $newHandoffZhCn
$wsZhCn.Range("K2").Value = $newHandbackDateZhCn

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = $newHandoffZhCn
$wsZhCn.Range("H3").Value = $newHandoffDateZhCn
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = $newHandoffZhCn
$wsZhCn.Range("K3").Value = $newHandbackDateZhCn

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/$oldFile1", [System.Type]::Missing, [System.Type]::Missing, $newFile1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ba657680638821218c5aaeacf1d40014b85cc67d/e2e/$oldFile1", [System.Type]::Missing, [System.Type]::Missing, $newFile1) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/$oldFile2", [System.Type]::Missing, [System.Type]::Missing, $newFile2) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ba657680638821218c5aaeacf1d40014b85cc67d/e2e/$oldFile2", [System.Type]::Missing, [System.Type]::Missing, $newFile2) | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = $newHandoffDeDe
$wsDeDe.Range("H2").Value = $newDate1
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = $newHandoffDeDe
$wsDeDe.Range("K2").Value = $newHandbackDateDeDe

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = $newHandoffDeDe
$wsDeDe.Range("H3").Value = $newDate1
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = $newHandoffDeDe
$wsDeDe.Range("K3").Value = $newHandbackDateDeDe

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/$oldFile1", [System.Type]::Missing, [System.Type]::Missing, $newFile1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d055dc71b68178f3101669a30fd340f2c93ceac2/e2e/$oldFile1", [System.Type]::Missing, [System.Type]::Missing, $newFile1) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68004ac1591356c7ffa5f9396679be1adf229d7f/e2e/$oldFile2", [System.Type]::Missing, [System.Type]::Missing, $newFile2) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d055dc71b68178f3101669a30fd340f2c93ceac2/e2e/$oldFile2", [System.Type]::Missing, [System.Type]::Missing, $newFile2) | Out-Null
